$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 0.3031743174317432
$ws.Range("D4").Value = 0.9850705070507049
$ws.Range("E4").Value = -0.9973477347734772
$ws.Range("G4").Value = -0.7628682868286828
$ws.Range("H4").Value = 0.3497749774977497
$ws.Range("I4").Value = -0.03313531353135313
$ws.Range("J4").Value = 0.2018001800180018
$ws.Range("K4").Value = -0.04021602160216022
$ws.Range("L4").Value = -0.05506150615061506
$ws.Range("M4").Value = -0.1050705070507051
$ws.Range("O4").Value = -0.04297629762976297
$ws.Range("P4").Value = -0.05143714371437143
$ws.Range("C5").Value = -0.1102190219021902
$ws.Range("D5").Value = -0.1154635463546354
$ws.Range("E5").Value = 0.1124632463246324
$ws.Range("G5").Value = 0.1175637563756376
$ws.Range("H5").Value = -0.2126852685268527
$ws.Range("I5").Value = -0.1426222622262226
$ws.Range("J5").Value = 0.08715271527152714
$ws.Range("K5").Value = 0.1484908490849085
$ws.Range("L5").Value = -0.01184518451845184
$ws.Range("M5").Value = 0.1243684368436844
$ws.Range("O5").Value = 0.09144914491449145
$ws.Range("P5").Value = 0.0318031803180318
$ws.Range("C6").Value = 0.1393579357935794
$ws.Range("D6").Value = 0.1216081608160816
$ws.Range("E6").Value = -0.1219321932193219
$ws.Range("G6").Value = -0.08153615361536153
$ws.Range("H6").Value = 0.03686768676867686
$ws.Range("I6").Value = -0.04927692769276927
$ws.Range("J6").Value = -0.04183618361836183
$ws.Range("K6").Value = 0.0255985598559856
$ws.Range("L6").Value = 0.0008880888088808879
$ws.Range("M6").Value = 0.1435463546354635
$ws.Range("O6").Value = 0.1348334833483348
$ws.Range("P6").Value = 0.09006900690069006
$ws.Range("C7").Value = 0.09209720972097209
$ws.Range("D7").Value = 0.1014341434143414
$ws.Range("E7").Value = 0.04727272727272726
$ws.Range("G7").Value = 0.4482568256825681
$ws.Range("H7").Value = 0.01466546654665466
$ws.Range("I7").Value = 0.007704770477047704
$ws.Range("J7").Value = -0.0297989798979898
$ws.Range("K7").Value = 0.2917251725172517
$ws.Range("L7").Value = 0.9896549654965496
$ws.Range("M7").Value = -0.1735613561356135
$ws.Range("O7").Value = 0.1419141914191419
$ws.Range("P7").Value = -0.09639363936393637
$ws.Range("C8").Value = 0.0568136813681368
$ws.Range("D8").Value = -0.06472247224722472
$ws.Range("E8").Value = 0.07237923792379237
$ws.Range("G8").Value = -0.3736933693369336
$ws.Range("H8").Value = 0.901182118211821
$ws.Range("I8").Value = 0.8634503450345034
$ws.Range("J8").Value = -0.04122412241224122
$ws.Range("K8").Value = 0.06573057305730572
$ws.Range("L8").Value = 0.008412841284128413
$ws.Range("M8").Value = 0.1332013201320132
$ws.Range("O8").Value = 0.03463546354635463
$ws.Range("P8").Value = -0.03805580558055805
$ws.Range("C9").Value = 0.8214101410141014
$ws.Range("D9").Value = 0.0248064806480648
$ws.Range("E9").Value = -0.03731173117311731
$ws.Range("G9").Value = -0.05136513651365136
$ws.Range("H9").Value = 0.02083408340834083
$ws.Range("I9").Value = 0.09438943894389439
$ws.Range("J9").Value = -0.04457245724572457
$ws.Range("K9").Value = -0.4868286828682867
$ws.Range("L9").Value = -0.03239123912391239
$ws.Range("M9").Value = -0.08802880288028801
$ws.Range("O9").Value = -0.01533753375337533
$ws.Range("P9").Value = -0.09437743774377437
$ws.Range("C10").Value = 0.2623702370237023
$ws.Range("D10").Value = -0.0488088808880888
$ws.Range("E10").Value = 0.06103810381038104
$ws.Range("G10").Value = 0.0454005400540054
$ws.Range("H10").Value = 0.0111011101110111
$ws.Range("I10").Value = -0.03294329432943294
$ws.Range("J10").Value = -0.06972697269726971
$ws.Range("K10").Value = 0.8430003000300028
$ws.Range("L10").Value = 0.07876387638763875
$ws.Range("M10").Value = 0.02774677467746775
$ws.Range("O10").Value = 0.03380738073807381
$ws.Range("P10").Value = 0.03648364836483647
$ws.Range("C11").Value = -0.09916591659165916
$ws.Range("D11").Value = 0.09646564656465645
$ws.Range("E11").Value = -0.1143114311431143
$ws.Range("G11").Value = -0.0444044404440444
$ws.Range("H11").Value = -0.1511071107110711
$ws.Range("I11").Value = -0.1291929192919292
$ws.Range("J11").Value = -0.01797779777977798
$ws.Range("K11").Value = -0.07641164116411642
$ws.Range("L11").Value = -0.1348454845484548
$ws.Range("M11").Value = 0.01624962496249625
$ws.Range("O11").Value = -0.03583558355835583
$ws.Range("P11").Value = -0.06186618661866185
$ws.Range("C12").Value = -0.1024182418241824
$ws.Range("D12").Value = -0.1870987098709871
$ws.Range("E12").Value = 0.1778217821782178
$ws.Range("G12").Value = 0.1774977497749775
$ws.Range("H12").Value = -0.1703930393039303
$ws.Range("I12").Value = -0.1041584158415841
$ws.Range("J12").Value = 0.05590159015901589
$ws.Range("K12").Value = 0.006492649264926492
$ws.Range("L12").Value = 0.02282628262826283
$ws.Range("M12").Value = -0.06024602460246024
$ws.Range("O12").Value = -0.01118511851185118
$ws.Range("P12").Value = 0.07853585358535851
$ws.Range("C13").Value = 0.2166456645664566
$ws.Range("D13").Value = 0.07679567956795678
$ws.Range("E13").Value = -0.07960396039603959
$ws.Range("G13").Value = -0.0516051605160516
$ws.Range("H13").Value = 0.0155055505550555
$ws.Range("I13").Value = 0.4652265226522652
$ws.Range("J13").Value = 0.09809780978097808
$ws.Range("K13").Value = -0.2077887788778878
$ws.Range("L13").Value = -0.01333333333333333
$ws.Range("M13").Value = -0.03635163516351635
$ws.Range("O13").Value = -0.02473447344734473
$ws.Range("P13").Value = -0.1393579357935794
$ws.Range("C14").Value = -0.1573957395739574
$ws.Range("D14").Value = -0.01382538253825382
$ws.Range("E14").Value = -0.01106510651065106
$ws.Range("G14").Value = -0.08930093009300928
$ws.Range("H14").Value = 0.02011401140114011
$ws.Range("I14").Value = -0.03866786678667866
$ws.Range("J14").Value = -0.08414041404140414
$ws.Range("K14").Value = -0.1655805580558056
$ws.Range("L14").Value = -0.1684728472847285
$ws.Range("M14").Value = 0.03710771077107711
$ws.Range("O14").Value = -0.03481548154815481
$ws.Range("P14").Value = 0.1923312331233123

Write-Host "Done updating cells."
